$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

$cell = $ws.Range("B1")
$cell.Borders.LineStyle = 1        # xlContinuous
$cell.Borders.Weight = 2           # xlThin
$cell.Font.Bold = $true
$cell.HorizontalAlignment = -4108  # xlCenter
$cell.VerticalAlignment = -4160    # xlTop

$cell.Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
